$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: extend the table with 4 new rows (21-24), cloning row 20s
#     full formatting (borders/fills/alignment) so the new rows visually
#     match the rest of the table.
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K24").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Step 2: fix up the IMPACT (column J) cell shading for the rows whose
#     threat level changes. Reuse the workbooks existing LOW/MEDIUM/HIGH
#     styles by copying the format from a cell that already carries the
#     right one, so fill/font exactly match the other LOW/MEDIUM/HIGH cells.

# NOTE: J19/J20 (new HIGH cells) are seeded from J18 *before* J18 itself is
# repainted to LOW below, so the copy source is still the original red style.
$ws.Range("J18").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J18").Copy()
$ws.Range("J20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("J7").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J2").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J7").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J7").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J2").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Step 3: keep column A as plain text (it holds "DD-MON-YY" strings,
#     not real dates) so re-assigning the values below does not get
#     auto-converted into date serials.
$ws.Range("A2:A24").NumberFormat = "@"

# --- Step 4: write the refreshed report data, row by row.

# Row 2: SM-444 vs Air Arabia Egypt E5-512 on 08-FEB-26
$ws.Range("A2").Value = "08-FEB-26"
$ws.Range("B2").Value = "SM-444"
$ws.Range("C2").Value = "Air Arabia Egypt E5-512"
$ws.Range("D2").Value = 612
$ws.Range("E2").Value = 647
$ws.Range("F2").Value = -35
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "LOW THREAT"
$ws.Range("K2").Value = "SAR"

# Row 3: SM-444 vs Nile Air NP-144 on 08-FEB-26
$ws.Range("A3").Value = "08-FEB-26"
$ws.Range("B3").Value = "SM-444"
$ws.Range("C3").Value = "Nile Air NP-144"
$ws.Range("D3").Value = 640
$ws.Range("E3").Value = 647
$ws.Range("F3").Value = -7
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

# Row 4: SM-444 vs Air Arabia Egypt E5-512 on 15-FEB-26
$ws.Range("A4").Value = "15-FEB-26"
$ws.Range("B4").Value = "SM-444"
$ws.Range("C4").Value = "Air Arabia Egypt E5-512"
$ws.Range("D4").Value = 612
$ws.Range("E4").Value = 647
$ws.Range("F4").Value = -35
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

# Row 5: SM-444 vs Nile Air NP-144 on 15-FEB-26
$ws.Range("A5").Value = "15-FEB-26"
$ws.Range("B5").Value = "SM-444"
$ws.Range("C5").Value = "Nile Air NP-144"
$ws.Range("D5").Value = 640
$ws.Range("E5").Value = 647
$ws.Range("F5").Value = -7
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

# Row 6: SM-448 vs Nile Air NP-144 on 19-FEB-26
$ws.Range("A6").Value = "19-FEB-26"
$ws.Range("B6").Value = "SM-448"
$ws.Range("C6").Value = "Nile Air NP-144"
$ws.Range("D6").Value = 571
$ws.Range("E6").Value = 591
$ws.Range("F6").Value = -20
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

# Row 7: SM-444 vs Air Arabia Egypt E5-512 on 20-FEB-26
$ws.Range("A7").Value = "20-FEB-26"
$ws.Range("B7").Value = "SM-444"
$ws.Range("C7").Value = "Air Arabia Egypt E5-512"
$ws.Range("D7").Value = 513
$ws.Range("E7").Value = 721
$ws.Range("F7").Value = -208
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K7").Value = "SAR"

# Row 8: SM-444 vs Air Arabia Egypt E5-512 on 22-FEB-26
$ws.Range("A8").Value = "22-FEB-26"
$ws.Range("B8").Value = "SM-444"
$ws.Range("C8").Value = "Air Arabia Egypt E5-512"
$ws.Range("D8").Value = 513
$ws.Range("E8").Value = 591
$ws.Range("F8").Value = -78
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "LOW THREAT"
$ws.Range("K8").Value = "SAR"

# Row 9: SM-444 vs Air Arabia Egypt E5-512 on 27-FEB-26
$ws.Range("A9").Value = "27-FEB-26"
$ws.Range("B9").Value = "SM-444"
$ws.Range("C9").Value = "Air Arabia Egypt E5-512"
$ws.Range("D9").Value = 513
$ws.Range("E9").Value = 883
$ws.Range("F9").Value = -370
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K9").Value = "SAR"

# Row 10: SM-444 vs Nile Air NP-144 on 01-MAR-26
$ws.Range("A10").Value = "01-MAR-26"
$ws.Range("B10").Value = "SM-444"
$ws.Range("C10").Value = "Nile Air NP-144"
$ws.Range("D10").Value = 971
$ws.Range("E10").Value = 1013
$ws.Range("F10").Value = -42
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "LOW THREAT"
$ws.Range("K10").Value = "SAR"

# Row 11: SM-444 vs Air Arabia Egypt E5-512 on 06-MAR-26
$ws.Range("A11").Value = "06-MAR-26"
$ws.Range("B11").Value = "SM-444"
$ws.Range("C11").Value = "Air Arabia Egypt E5-512"
$ws.Range("D11").Value = 934
$ws.Range("E11").Value = 1159
$ws.Range("F11").Value = -225
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K11").Value = "SAR"

# Row 12: SM-444 vs Air Arabia Egypt E5-512 on 08-MAR-26
$ws.Range("A12").Value = "08-MAR-26"
$ws.Range("B12").Value = "SM-444"
$ws.Range("C12").Value = "Air Arabia Egypt E5-512"
$ws.Range("D12").Value = 934
$ws.Range("E12").Value = 1013
$ws.Range("F12").Value = -79
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "LOW THREAT"
$ws.Range("K12").Value = "SAR"

# Row 13: SM-444 vs Air Arabia Egypt E5-512 on 15-MAR-26
$ws.Range("A13").Value = "15-MAR-26"
$ws.Range("B13").Value = "SM-444"
$ws.Range("C13").Value = "Air Arabia Egypt E5-512"
$ws.Range("D13").Value = 1236
$ws.Range("E13").Value = 1306
$ws.Range("F13").Value = -70
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "LOW THREAT"
$ws.Range("K13").Value = "SAR"

# Row 14: SM-444 vs Air Arabia Egypt E5-512 on 20-MAR-26
$ws.Range("A14").Value = "20-MAR-26"
$ws.Range("B14").Value = "SM-444"
$ws.Range("C14").Value = "Air Arabia Egypt E5-512"
$ws.Range("D14").Value = 1236
$ws.Range("E14").Value = 1306
$ws.Range("F14").Value = -70
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "LOW THREAT"
$ws.Range("K14").Value = "SAR"

# Row 15: SM-444 vs Air Arabia Egypt E5-512 on 27-MAR-26
$ws.Range("A15").Value = "27-MAR-26"
$ws.Range("B15").Value = "SM-444"
$ws.Range("C15").Value = "Air Arabia Egypt E5-512"
$ws.Range("D15").Value = 513
$ws.Range("E15").Value = 721
$ws.Range("F15").Value = -208
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = 30
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K15").Value = "SAR"

# Row 16: SM-444 vs Air Arabia Egypt E5-512 on 30-MAR-26
$ws.Range("A16").Value = "30-MAR-26"
$ws.Range("B16").Value = "SM-444"
$ws.Range("C16").Value = "Air Arabia Egypt E5-512"
$ws.Range("D16").Value = 507
$ws.Range("E16").Value = 648
$ws.Range("F16").Value = -141
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "LOW THREAT"
$ws.Range("K16").Value = "SAR"

# Row 17: SM-444 vs Nile Air NP-144 on 01-APR-26
$ws.Range("A17").Value = "01-APR-26"
$ws.Range("B17").Value = "SM-444"
$ws.Range("C17").Value = "Nile Air NP-144"
$ws.Range("D17").Value = 571
$ws.Range("E17").Value = 591
$ws.Range("F17").Value = -20
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "LOW THREAT"
$ws.Range("K17").Value = "SAR"

# Row 18: SM-444 vs Air Arabia Egypt E5-512 on 02-APR-26
$ws.Range("A18").Value = "02-APR-26"
$ws.Range("B18").Value = "SM-444"
$ws.Range("C18").Value = "Air Arabia Egypt E5-512"
$ws.Range("D18").Value = 507
$ws.Range("E18").Value = 564
$ws.Range("F18").Value = -57
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "LOW THREAT"
$ws.Range("K18").Value = "SAR"

# Row 19: SM-444 vs Air Arabia Egypt E5-514 on 16-MAY-26
$ws.Range("A19").Value = "16-MAY-26"
$ws.Range("B19").Value = "SM-444"
$ws.Range("C19").Value = "Air Arabia Egypt E5-514"
$ws.Range("D19").Value = 934
$ws.Range("E19").Value = 1501
$ws.Range("F19").Value = -567
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 30
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K19").Value = "SAR"

# Row 20: SM-444 vs Air Arabia Egypt E5-512 on 18-MAY-26
$ws.Range("A20").Value = "18-MAY-26"
$ws.Range("B20").Value = "SM-444"
$ws.Range("C20").Value = "Air Arabia Egypt E5-512"
$ws.Range("D20").Value = 928
$ws.Range("E20").Value = 1501
$ws.Range("F20").Value = -573
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K20").Value = "SAR"

# Row 21: SM-444 vs Nile Air NP-144 on 20-MAY-26
$ws.Range("A21").Value = "20-MAY-26"
$ws.Range("B21").Value = "SM-444"
$ws.Range("C21").Value = "Nile Air NP-144"
$ws.Range("D21").Value = 1271
$ws.Range("E21").Value = 1501
$ws.Range("F21").Value = -230
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K21").Value = "SAR"

# Row 22: SM-444 vs Air Arabia Egypt E5-512 on 21-MAY-26
$ws.Range("A22").Value = "21-MAY-26"
$ws.Range("B22").Value = "SM-444"
$ws.Range("C22").Value = "Air Arabia Egypt E5-512"
$ws.Range("D22").Value = 1231
$ws.Range("E22").Value = 1501
$ws.Range("F22").Value = -270
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = 30
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K22").Value = "SAR"

# Row 23: SM-444 vs Air Arabia Egypt E5-514 on 23-MAY-26
$ws.Range("A23").Value = "23-MAY-26"
$ws.Range("B23").Value = "SM-444"
$ws.Range("C23").Value = "Air Arabia Egypt E5-514"
$ws.Range("D23").Value = 1236
$ws.Range("E23").Value = 1501
$ws.Range("F23").Value = -265
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K23").Value = "SAR"

# Row 24: SM-444 vs Air Arabia Egypt E5-512 on 25-MAY-26
$ws.Range("A24").Value = "25-MAY-26"
$ws.Range("B24").Value = "SM-444"
$ws.Range("C24").Value = "Air Arabia Egypt E5-512"
$ws.Range("D24").Value = 1231
$ws.Range("E24").Value = 1501
$ws.Range("F24").Value = -270
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K24").Value = "SAR"
